# Adds 4 new CRESS virus reference rows (Redondoviridae/Smacoviridae and
# related lineages) to Sheet1, as described in the commit:
#   "Updated to include references for Redondoviridae, Smacoviridae,
#    among other lineages"
#
# Columns: A=sequenceID  B=name  C=full_name  D=family  E=genus  F=clade
#          G=host range  H=isolation_source

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Write the new cell values. The order in which *new* distinct strings
#    are first assigned controls the order they are appended to the
#    shared-string table, so we deliberately set them in the same order
#    the source workbook's sharedStrings.xml grew in.
# ---------------------------------------------------------------------

# -- Column A (sequenceID) for the 4 new rows, top to bottom --
$ws.Cells.Item(50, 1).Value = "MK059763"
$ws.Cells.Item(51, 1).Value = "MF669476"
$ws.Cells.Item(52, 1).Value = "KJ888053"
$ws.Cells.Item(53, 1).Value = "MK059757"

# -- family / genus / clade block, entered bottom row first --
$ws.Cells.Item(53, 5).Value = "Brisavirus"
$ws.Cells.Item(53, 4).Value = "Brisavirus"
$ws.Cells.Item(53, 6).Value = "Brisavirus"

$ws.Cells.Item(52, 5).Value = "Gemininvirus"
$ws.Cells.Item(52, 4).Value = "Gemininviridae"
$ws.Cells.Item(52, 6).Value = "Gemininvirus"

$ws.Cells.Item(51, 5).Value = "Smacovirus"
$ws.Cells.Item(51, 4).Value = "Smacovirus"
$ws.Cells.Item(51, 6).Value = "Smacovirus"

$ws.Cells.Item(50, 5).Value = "Vientovirus"
$ws.Cells.Item(50, 4).Value = "Vientovirus"
$ws.Cells.Item(50, 6).Value = "Vientovirus"

# -- Column B (name) for the 4 new rows, top to bottom --
$ws.Cells.Item(50, 2).Value = "HLVV"
$ws.Cells.Item(51, 2).Value = "CBSV"
$ws.Cells.Item(52, 2).Value = "EACMKV"
$ws.Cells.Item(53, 2).Value = "L-BrV-RC"

# -- Column C (full_name), entered bottom row first --
$ws.Cells.Item(53, 3).Value = "Human lung-associated brisavirus RC"
$ws.Cells.Item(52, 3).Value = "East African cassava mosaic Kenya virus"
$ws.Cells.Item(51, 3).Value = "Cattle blood-associated circovirus-like virus "
$ws.Cells.Item(50, 3).Value = "Human lung-associated vientovirus FB"

# -- Columns G/H (host range / isolation_source) reuse the existing
#    "Unknown" shared string already present in the workbook --
$ws.Cells.Item(50, 7).Value = "Unknown"
$ws.Cells.Item(50, 8).Value = "Unknown"
$ws.Cells.Item(51, 7).Value = "Unknown"
$ws.Cells.Item(51, 8).Value = "Unknown"
$ws.Cells.Item(52, 7).Value = "Unknown"
$ws.Cells.Item(52, 8).Value = "Unknown"
$ws.Cells.Item(53, 7).Value = "Unknown"
$ws.Cells.Item(53, 8).Value = "Unknown"

# ---------------------------------------------------------------------
# 2. Copy cell formatting down from the last pre-existing data row (49)
#    so the new rows carry the same visual style as the rest of the
#    table: shaded A:C, plain D:F, shaded G:H.
# ---------------------------------------------------------------------
$ws.Range("A49:C49").Copy() | Out-Null
$ws.Range("A50:C53").PasteSpecial(-4122) | Out-Null

$ws.Range("G49:H49").Copy() | Out-Null
$ws.Range("G50:H53").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Refresh the worksheet selection / view to cover the new range.
# ---------------------------------------------------------------------
$ws.Range("A1:H53").Select() | Out-Null

Write-Output "Added rows 50:53 to Sheet1"
